$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Append a new row (11) to the Login data-driven-test table
$ws.Range("A11").Value = "test@test.com"
$ws.Range("B11").Value = "test123"

$ws.Range("A11").Select()
